$d = $word.ActiveDocument

# 1. "Azure Account" -> "Azure Account Subscription"
$d.Content.Find.Execute("Azure Account", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Azure Account Subscription", 2) | Out-Null

# 2. Azure CLI 2.0 bullet rewrite
$d.Content.Find.Execute(
    "Azure CLI 2.0 (in your local system, if possible). If you don" + [char]0x2019 + "t want to you can use the one provided by Azure Portal to do most of the things",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Azure CLI 2.0 (in your local system, if possible)- Either portal Azure CLI or download the Azure CLI 2.0.", 2) | Out-Null

# 3. Kubectl.exe bullet: drop ". This is needed" before "to see Kubernetes Dashboard"
$d.Content.Find.Execute(
    " command). This is needed to see Kubernetes Dashboard",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " command) to see Kubernetes Dashboard", 2) | Out-Null

# 4. "nt to join our existing cluster" -> "nt to join the existing cluster"
$d.Content.Find.Execute("nt to join our existing cluster", $true, $false, $false, $false, $false,
                         $true, 1, $false, "nt to join the existing cluster", 2) | Out-Null

# 5. "you will see a cluster successfully" -> "you will note a cluster successfully"
$d.Content.Find.Execute("you will see a cluster successfully", $true, $false, $false, $false, $false,
                         $true, 1, $false, "you will note a cluster successfully", 2) | Out-Null

# 6. "After this much , write the following command" -> "After this, write the following command"
$d.Content.Find.Execute("After this much , write the following command", $true, $false, $false, $false, $false,
                         $true, 1, $false, "After this, write the following command", 2) | Out-Null

# 7. "cluster wont be setup properly" -> "cluster won't be setup properly" (curly apostrophe)
$d.Content.Find.Execute("otherwise cluster wont be setup properly", $true, $false, $false, $false, $false,
                         $true, 1, $false, ("otherwise cluster won" + [char]0x2019 + "t be setup properly"), 2) | Out-Null

# 8. "Kubernetes Dashboard after the setup, ill recommend to see" -> "... its recommended to see"
$d.Content.Find.Execute("Kubernetes Dashboard after the setup, ill recommend to see",
                         $true, $false, $false, $false, $false, $true, 1, $false,
                         "Kubernetes Dashboard after the setup, its recommended to see", 2) | Out-Null

# 9. "In order to connect ... kindly copy and paste this " rewrite
$d.Content.Find.Execute(
    "In order to connect to the existing cluster which is already deployed by us from your local or azure portal , kindly copy and paste this ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In order to connect to the existing cluster from your local or azure portal which is already up and running , copy and paste the ", 2) | Out-Null

# 10. "Once you this and see th" -> "Once you execute this and see th"
$d.Content.Find.Execute("Once you this and see th", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Once you execute this and see th", 2) | Out-Null

Write-Output "simple replacements done"
